# Applies the Convocation.docx text edits described by the commit diff.
# Six distinct text substitutions, each occurring twice in the document
# (the document repeats the same letter body twice).

$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

function Replace-AllText($findText, $replaceText) {
    $d.Content.Find.Execute(
        $findText,    # FindText
        $true,        # MatchCase
        $false,       # MatchWholeWord
        $false,       # MatchWildcards
        $false,       # MatchSoundsLike
        $false,       # MatchAllWordForms
        $true,        # Forward
        1,            # Wrap (wdFindContinue)
        $false,       # Format
        $replaceText, # ReplaceWith
        $wdReplaceAll # Replace
    ) | Out-Null
}

# --- "محمد" is ambiguous: it also occurs as a substring inside the
# unrelated signature line "قدور بن دهمة محمد الأمين", which must stay
# untouched. Scope that replacement to only the paragraphs that also
# contain "جبور" (the salutation line), which is where the standalone
# "محمد" run actually lives, right next to it.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*جبور*") {
        $p.Range.Find.Execute(
            "محمد",
            $true,
            $true,
            $false,
            $false,
            $false,
            $true,
            1,
            $false,
            "ddddddddd",
            $wdReplaceAll
        ) | Out-Null
    }
}

# --- Remaining five substitutions are unambiguous document-wide.
Replace-AllText "جبور " "dddddddd "
Replace-AllText "دوار أهل البرادعي" "ddddddddddddd"
Replace-AllText "واد الصباح" "عين تموشنت"
Replace-AllText "نسخة من بطاقة التعريف الوطنية + شهادة عدم الانتساب للضمان الإجتماعي للأجراء ( NON AFFILIATION CNAS )" "شهادة عدم الإنتساب للضمان الإجتماعي ( CNAS )"
Replace-AllText "من أجل تسوية وضعية المنحة المالية للمعوقين" "Mot du dass"

Write-Output "Done"
